# Reorder the monthly rows (A2:C49) so that, within each year, the
# October/November/December rows move to the front of that year's block
# (ahead of January..September), while keeping each row's own A/B/C
# values together and unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:C49")
$arr = $rng.Value2
$n = $arr.GetLength(0)

# Distinct years, in order of first appearance.
$years = @()
for ($i = 1; $i -le $n; $i++) {
    $label = $arr[$i, 1]
    $year = $label.Substring(0, 4)
    if ($years -notcontains $year) {
        $years += $year
    }
}

$newA = @()
$newB = @()
$newC = @()

foreach ($year in $years) {
    # Oct/Nov/Dec rows for this year first, in their original relative order.
    for ($i = 1; $i -le $n; $i++) {
        $label = $arr[$i, 1]
        if ($label.Substring(0, 4) -eq $year) {
            $month = $label.Substring(5, 2)
            if ($month -eq "10" -or $month -eq "11" -or $month -eq "12") {
                $newA += $arr[$i, 1]
                $newB += $arr[$i, 2]
                $newC += $arr[$i, 3]
            }
        }
    }
    # Then Jan..Sep rows for this year, in their original relative order.
    for ($i = 1; $i -le $n; $i++) {
        $label = $arr[$i, 1]
        if ($label.Substring(0, 4) -eq $year) {
            $month = $label.Substring(5, 2)
            if ($month -ne "10" -and $month -ne "11" -and $month -ne "12") {
                $newA += $arr[$i, 1]
                $newB += $arr[$i, 2]
                $newC += $arr[$i, 3]
            }
        }
    }
}

for ($i = 1; $i -le $n; $i++) {
    $arr[$i, 1] = $newA[$i - 1]
    $arr[$i, 2] = $newB[$i - 1]
    $arr[$i, 3] = $newC[$i - 1]
}

$rng.Value2 = $arr
